# Update the daily crypto price/volume snapshot (columns D and E) with the
# latest scraped figures from the GitHub Actions run.
#
# Every cell in columns D/E is stored as text in the source workbook (prices
# such as "2.819.00" use '.' as both thousands- and decimal-separator, so a
# naive numeric assignment would corrupt them). Excel's COM layer still tries
# to auto-coerce plain decimal-looking strings (e.g. "187.17", "1.00") into
# numbers, so for those we briefly force the cell to Text format before the
# write and restore the default ("Normal") style afterwards so the cell
# formatting itself is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Cell,
        [string]$Text
    )

    # Strings that look like a single plain decimal number (e.g. "187.17",
    # "1.00") get silently reinterpreted as a Double by Excel's COM layer.
    # Anything else (multiple '.' as in "74.862.41", or values containing
    # '%'/spaces) is never auto-coerced and can be assigned directly.
    $looksNumeric = $Text -match '^-?\d+(\.\d+)?$'

    $range = $ws.Range($Cell)
    if ($looksNumeric) {
        $range.NumberFormat = "@"
        $range.Value = $Text
        $range.Style = "Normal"
    }
    else {
        $range.Value = $Text
    }
}

# --- Row 2 (Bitcoin) ---
Set-TextValue "D2" "74.862.41"
Set-TextValue "E2" "  +0.79%  "

# --- Row 3 (Ethereum) ---
Set-TextValue "D3" "2.816.74"
Set-TextValue "E3" "  +6.94%  "

# --- Row 4 (TetherUSD) ---
Set-TextValue "E4" "  +0.03%  "

# --- Row 5 (Solana) ---
Set-TextValue "D5" "187.17"
Set-TextValue "E5" "  +0.69%  "

# --- Row 6 (BNB) ---
Set-TextValue "D6" "594.16"
Set-TextValue "E6" "  +1.94%  "

# --- Row 7 (USDC) ---
Set-TextValue "E7" "  +0.04%  "

# --- Row 8 (XRP) ---
Set-TextValue "D8" "0.548"
Set-TextValue "E8" "  +2.63%  "

# --- Row 9 (Dogecoin) ---
Set-TextValue "E9" "  -4.86%  "

# --- Row 10 (LidoStakedEther) ---
Set-TextValue "D10" "2.814.10"
Set-TextValue "E10" "  +6.88%  "

# --- Row 11 (TRON) ---
Set-TextValue "D11" "0.160"
Set-TextValue "E11" "  -1.35%  "

# --- Row 12 (Cardano) ---
Set-TextValue "E12" "  +3.54%  "

# --- Row 13 (Toncoin) ---
Set-TextValue "E13" "  +2.80%  "

# --- Row 14 (WrappedliquidstakedEther2.0) ---
Set-TextValue "D14" "3.337.88"
Set-TextValue "E14" "  +7.14%  "

# --- Row 15 (WrappedBTC) ---
Set-TextValue "D15" "74.776.53"
Set-TextValue "E15" "  +0.93%  "

# --- Row 16 (ShibaInu) ---
Set-TextValue "E16" "  -1.42%  "

# --- Row 17 (Avalanche) ---
Set-TextValue "D17" "26.78"
Set-TextValue "E17" "  +1.82%  "

# --- Row 18 (WrappedEther) ---
Set-TextValue "D18" "2.815.18"
Set-TextValue "E18" "  +6.34%  "

# --- Row 19 (Uniswap) ---
Set-TextValue "D19" "8.95"
Set-TextValue "E19" "  -1.64%  "

# --- Row 20 (Chainlink) ---
Set-TextValue "E20" "  +3.97%  "

# --- Row 21 (BitcoinCash) ---
Set-TextValue "D21" "377.36"
Set-TextValue "E21" "  +1.30%  "

# --- Row 22 (SuiNetwork) ---
Set-TextValue "E22" "  -1.82%  "

# --- Row 23 (Polkadot) ---
Set-TextValue "E23" "  -0.71%  "

# --- Row 24 (Dai) ---
Set-TextValue "E24" "  -0.06%  "

# --- Row 25 (Litecoin) ---
Set-TextValue "D25" "70.89"
Set-TextValue "E25" "  +1.08%  "

# --- Row 26 (Aptos) ---
Set-TextValue "D26" "9.87"
Set-TextValue "E26" "  +5.21%  "

# --- Row 27 (WrappedeETH) ---
Set-TextValue "E27" "  +7.15%  "

# --- Row 28 (NEARProtocol) ---
Set-TextValue "E28" "  -0.25%  "

# --- Row 29 (PEPE) ---
Set-TextValue "E29" "  +9.05%  "

# --- Row 30 (Binance-PegBSC-USD) ---
Set-TextValue "D30" "0.997"
Set-TextValue "E30" "  -0.52%  "

# --- Row 31 (Bittensor) ---
Set-TextValue "D31" "517.13"
Set-TextValue "E31" "  -1.87%  "

# --- Row 32 (Fetch.AI) ---
Set-TextValue "E32" "  -0.43%  "

# --- Row 33 (InternetComputer(DFINITY)) ---
Set-TextValue "E33" "  +0.07%  "

# --- Row 34 (PancakeSwap) ---
Set-TextValue "E34" "  +2.52%  "

# --- Row 35 (FirstDigitalUSD) ---
Set-TextValue "D35" "1.00"
Set-TextValue "E35" "  +0.03%  "

# --- Row 36 (Monero) ---
Set-TextValue "D36" "163.53"
Set-TextValue "E36" "  +0.16%  "

# --- Row 37 (EthereumClassic) ---
Set-TextValue "E37" "  +3.87%  "

# --- Row 38 (Kaspa) ---
Set-TextValue "E38" "  -1.17%  "

# --- Row 39 (WhiteBITCoin) ---
Set-TextValue "E39" "  +0.45%  "

# --- Row 40 (Aave) ---
Set-TextValue "D40" "186.95"
Set-TextValue "E40" "  +16.41%  "

# --- Row 41 (USDe) ---
Set-TextValue "E41" "  -0.02%  "

# --- Row 42 (PolygonEcosystemToken) ---
Set-TextValue "D42" "0.341"

# --- Row 43 (RenderToken) ---
Set-TextValue "D43" "4.99"
Set-TextValue "E43" "  +1.65%  "

# --- Row 44 (Stacks) ---
Set-TextValue "E44" "  -0.42%  "

# --- Row 45 (ImmutableX) ---
Set-TextValue "E45" "  +1.79%  "

# --- Row 46 (OKB) ---
Set-TextValue "D46" "39.93"
Set-TextValue "E46" "  +2.41%  "

# --- Row 47 (Cronos) ---
Set-TextValue "E47" "  +0.21%  "

# --- Row 48 (dogwifhat) ---
Set-TextValue "E48" "  -2.52%  "

# --- Row 49 (ARBITRUM) ---
Set-TextValue "E49" "  +8.98%  "

# --- Row 50 (Filecoin) ---
Set-TextValue "D50" "3.71"
Set-TextValue "E50" "  +2.25%  "

# --- Row 51 (Mantle) ---
Set-TextValue "E51" "  +7.96%  "
